# Run prelim process with final data set - set seed for imputation for reproducibility
# Update Paper_ID (column A) values for rows 9-140 to reflect the re-run with the
# final data set (IDs shift upward because additional papers were incorporated
# earlier in the pipeline).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of contiguous row ranges (start,end) -> new Paper_ID value
$updates = @(
    @{ Start = 9;   End = 12;  New = 102 },
    @{ Start = 13;  End = 18;  New = 107 },
    @{ Start = 19;  End = 33;  New = 118 },
    @{ Start = 34;  End = 45;  New = 137 },
    @{ Start = 46;  End = 51;  New = 18  },
    @{ Start = 52;  End = 68;  New = 91  },
    @{ Start = 69;  End = 74;  New = 138 },
    @{ Start = 75;  End = 76;  New = 140 },
    @{ Start = 77;  End = 81;  New = 118 },
    @{ Start = 82;  End = 116; New = 128 },
    @{ Start = 117; End = 122; New = 129 },
    @{ Start = 123; End = 140; New = 150 }
)

foreach ($u in $updates) {
    $rangeAddress = "A" + $u.Start + ":A" + $u.End
    $ws.Range($rangeAddress).Value = $u.New
}
